# Natmi following Dr Hou advice
# Update existing rows 2-4 with recomputed NATMI edge-weight stats, and
# split the previous single "sCs" target-cluster row into three rows:
# M1, M2 and sCs (row 4 becomes M1, new rows 5 and 6 become M2 and sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (Target cluster: ECs) ----
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3616933333333334
$ws.Range("H2").Value = 1.08508
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.784940333333333
$ws.Range("N2").Value = 5.354821
$ws.Range("O2").Value = 0.2226931997744963
$ws.Range("P2").Value = 0.2226931997744964
$ws.Range("Q2").Value = 0.6456010189644446
$ws.Range("R2").Value = 5.810409170680001
$ws.Range("S2").Value = 0.2226931997744963
$ws.Range("T2").Value = 0.2226931997744964

# ---- Row 3 (Target cluster: FAPs) ----
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3616933333333334
$ws.Range("H3").Value = 1.08508
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.575715666666667
$ws.Range("N3").Value = 10.727147
$ws.Range("O3").Value = 0.4461143873682032
$ws.Range("P3").Value = 0.4461143873682032
$ws.Range("Q3").Value = 1.293312518528889
$ws.Range("R3").Value = 11.63981266676
$ws.Range("S3").Value = 0.4461143873682032
$ws.Range("T3").Value = 0.4461143873682032

# ---- Row 4 (Target cluster: was sCs, now M1) ----
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3616933333333334
$ws.Range("H4").Value = 1.08508
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02040433333333333
$ws.Range("N4").Value = 0.061213
$ws.Range("O4").Value = 0.002545690852746758
$ws.Range("P4").Value = 0.002545690852746758
$ws.Range("Q4").Value = 0.007380111337777779
$ws.Range("R4").Value = 0.06642100204
$ws.Range("S4").Value = 0.002545690852746758
$ws.Range("T4").Value = 0.002545690852746758

# ---- Row 5 (new, Target cluster: M2) ----
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Tgfa"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3616933333333334
$ws.Range("H5").Value = 1.08508
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.020115
$ws.Range("N5").Value = 0.060345
$ws.Range("O5").Value = 0.002509592970594532
$ws.Range("P5").Value = 0.002509592970594532
$ws.Range("Q5").Value = 0.007275461400000001
$ws.Range("R5").Value = 0.06547915260000001
$ws.Range("S5").Value = 0.002509592970594532
$ws.Range("T5").Value = 0.002509592970594532

# ---- Row 6 (new, Target cluster: sCs) ----
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Tgfa"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3616933333333334
$ws.Range("H6").Value = 1.08508
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.614068666666667
$ws.Range("N6").Value = 7.842206
$ws.Range("O6").Value = 0.3261371290339591
$ws.Range("P6").Value = 0.3261371290339591
$ws.Range("Q6").Value = 0.945491209608889
$ws.Range("R6").Value = 8.509420886480001
$ws.Range("S6").Value = 0.3261371290339591
$ws.Range("T6").Value = 0.3261371290339591
